$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wishlist row (row 15): Aeron Chair
# Column order matches the order new shared strings were appended in the
# original commit: B (Image URL), A (Name), C (Link URL), D (Price)
$ws.Range("B15").Value = "https://images.hermanmiller.group/m/49e8a45a115fae9a/W-HM_AER_61524.png?trim=auto&trim-sd=1&blend-mode=darken&blend=fafafa&bg=fafafa&auto=format&w=1000&q=70&h=1000"
$ws.Range("A15").Value = "Aeron Chair"
$ws.Range("C15").Value = "https://store.hermanmiller.com/office-chairs-aeron/aeron-chair/100249179.html?lang=en_US&sku=100249179"
$ws.Range("D15").Value = 1700

$ws.Range("D16").Select()
